{"js": "// Revert \"Agregado de URL de repositorio\":\n// 1) Remove the \" (https://github.com/leomalevo/TiendaEvertec )\" parenthetical\n//    (together with its hyperlink) that was appended after\n//    \"Baje del repositorio el c\u00f3digo fuente de la soluci\u00f3n web\".\n// 2) Remove the \"TiendaEvertec.sln \" mention from the step that talks about\n//    opening the solution in Visual Studio.\n\nconst body = context.document.body;\n\n// --- Change 1: drop the repo URL parenthetical + hyperlink -----------------\nconst urlRanges = body.search(\" (https://github.com/leomalevo/TiendaEvertec )\", {\n  matchCase: true\n});\nurlRanges.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < urlRanges.items.length; i++) {\n  urlRanges.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 2: drop the \"TiendaEvertec.sln \" mention ------------------------\nconst slnRanges = body.search(\"TiendaEvertec.sln \", { matchCase: true });\nslnRanges.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < slnRanges.items.length; i++) {\n  slnRanges.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Revert \"Agregado de URL de repositorio\":\n# 1) Remove the \" (https://github.com/leomalevo/TiendaEvertec )\" parenthetical\n#    (together with its hyperlink) that was appended after\n#    \"Baje del repositorio el c\u00f3digo fuente de la soluci\u00f3n web\".\n# 2) Remove the \"TiendaEvertec.sln \" mention from the step that talks about\n#    opening the solution in Visual Studio.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: drop the repo URL parenthetical + hyperlink -----------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\" (https://github.com/leomalevo/TiendaEvertec )\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# --- Change 2: drop the \"TiendaEvertec.sln \" mention ------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"TiendaEvertec.sln \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n"}
